# Auto-generated edit script: update Leve profit-calculation columns (H-N)
# across multiple worksheets per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H127").Value = 938
$ws.Range("J127").Value = 940.6667
$ws.Range("L127").Value = 2822.0001
$ws.Range("N127").Value = -12742.0001
$ws.Range("H137").Value = 2108.7693
$ws.Range("I137").Value = 1162.75
$ws.Range("J137").Value = 2280.7727
$ws.Range("K137").Value = 3488.25
$ws.Range("L137").Value = 6842.3181
$ws.Range("M137").Value = -938.25
$ws.Range("N137").Value = -11942.3181
$ws.Range("H138").Value = 2169.8372
$ws.Range("I138").Value = 2080.261
$ws.Range("J138").Value = 2272.85
$ws.Range("K138").Value = 6240.782999999999
$ws.Range("L138").Value = 6818.549999999999
$ws.Range("M138").Value = -1100.782999999999
$ws.Range("N138").Value = -17098.55
$ws.Range("H141").Value = 1335792.5
$ws.Range("I141").Value = 1649155.4
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 4947466.199999999
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -4942286.199999999
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2339.7402
$ws.Range("I32").Value = 1670.4697
$ws.Range("J32").Value = 6355.364
$ws.Range("K32").Value = 1670.4697
$ws.Range("L32").Value = 6355.364
$ws.Range("M32").Value = -1383.4697
$ws.Range("N32").Value = -6929.364
$ws.Range("H45").Value = 3106146
$ws.Range("I45").Value = 7502502.5
$ws.Range("K45").Value = 7502502.5
$ws.Range("M45").Value = -7502125.5
$ws.Range("H61").Value = 4473.4614
$ws.Range("I61").Value = 2739.4
$ws.Range("J61").Value = 5557.25
$ws.Range("K61").Value = 2739.4
$ws.Range("L61").Value = 5557.25
$ws.Range("M61").Value = -2527.4
$ws.Range("N61").Value = -5981.25
$ws.Range("H74").Value = 1843.0834
$ws.Range("I74").Value = 629
$ws.Range("K74").Value = 629
$ws.Range("M74").Value = 245
$ws.Range("H77").Value = 1843.0834
$ws.Range("I77").Value = 629
$ws.Range("K77").Value = 3145
$ws.Range("M77").Value = 1223
$ws.Range("H97").Value = 1867.8
$ws.Range("I97").Value = 1786.9286
$ws.Range("K97").Value = 1786.9286
$ws.Range("M97").Value = -1290.9286
$ws.Range("H122").Value = 48945.125
$ws.Range("I122").Value = 96370.25
$ws.Range("K122").Value = 289110.75
$ws.Range("M122").Value = -286660.75
$ws.Range("H136").Value = 4473.4614
$ws.Range("I136").Value = 2739.4
$ws.Range("J136").Value = 5557.25
$ws.Range("K136").Value = 8218.200000000001
$ws.Range("L136").Value = 16671.75
$ws.Range("M136").Value = -5668.200000000001
$ws.Range("N136").Value = -21771.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11593.909
$ws.Range("I134").Value = 12353.3
$ws.Range("K134").Value = 37059.89999999999
$ws.Range("M134").Value = -34524.89999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1376.6
$ws.Range("I31").Value = 871.9286
$ws.Range("J31").Value = 1713.0476
$ws.Range("K31").Value = 871.9286
$ws.Range("L31").Value = 1713.0476
$ws.Range("M31").Value = -576.9286
$ws.Range("N31").Value = -2303.0476
$ws.Range("H34").Value = 1376.6
$ws.Range("I34").Value = 871.9286
$ws.Range("J34").Value = 1713.0476
$ws.Range("K34").Value = 871.9286
$ws.Range("L34").Value = 1713.0476
$ws.Range("M34").Value = -669.9286
$ws.Range("N34").Value = -2117.0476
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H63").Value = 36270
$ws.Range("J63").Value = 36270
$ws.Range("L63").Value = 36270
$ws.Range("N63").Value = -37642
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H66").Value = 36270
$ws.Range("J66").Value = 36270
$ws.Range("L66").Value = 108810
$ws.Range("N66").Value = -115674
$ws.Range("H99").Value = 2837.5
$ws.Range("I99").Value = 2550
$ws.Range("J99").Value = 3125
$ws.Range("K99").Value = 2550
$ws.Range("L99").Value = 3125
$ws.Range("M99").Value = -1052
$ws.Range("N99").Value = -6121
$ws.Range("H122").Value = 1055
$ws.Range("I122").Value = 1055
$ws.Range("K122").Value = 3165
$ws.Range("M122").Value = -715
$ws.Range("H126").Value = 2837.5
$ws.Range("I126").Value = 2550
$ws.Range("J126").Value = 3125
$ws.Range("K126").Value = 7650
$ws.Range("L126").Value = 9375
$ws.Range("M126").Value = -5180
$ws.Range("N126").Value = -14315
$ws.Range("H132").Value = 2758.5
$ws.Range("I132").Value = 1890.3
$ws.Range("K132").Value = 5670.9
$ws.Range("M132").Value = -3140.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 354.45
$ws.Range("I5").Value = 337.125
$ws.Range("J5").Value = 366
$ws.Range("K5").Value = 1011.375
$ws.Range("L5").Value = 1098
$ws.Range("M5").Value = -899.375
$ws.Range("N5").Value = -1322
$ws.Range("H12").Value = 53.77778
$ws.Range("J12").Value = 112.14286
$ws.Range("L12").Value = 336.42858
$ws.Range("N12").Value = -682.42858
$ws.Range("H38").Value = 260.8
$ws.Range("I38").Value = 61.4
$ws.Range("J38").Value = 460.2
$ws.Range("K38").Value = 184.2
$ws.Range("L38").Value = 1380.6
$ws.Range("M38").Value = 162.8
$ws.Range("N38").Value = -2074.6
$ws.Range("H68").Value = 1625.0426
$ws.Range("I68").Value = 779.375
$ws.Range("J68").Value = 1798.5128
$ws.Range("K68").Value = 2338.125
$ws.Range("L68").Value = 5395.538399999999
$ws.Range("M68").Value = -1527.125
$ws.Range("N68").Value = -7017.538399999999
$ws.Range("H71").Value = 1625.0426
$ws.Range("I71").Value = 779.375
$ws.Range("J71").Value = 1798.5128
$ws.Range("K71").Value = 7014.375
$ws.Range("L71").Value = 16186.6152
$ws.Range("M71").Value = -2958.375
$ws.Range("N71").Value = -24298.6152
$ws.Range("H107").Value = 1427.4736
$ws.Range("I107").Value = 1271.4166
$ws.Range("J107").Value = 1695
$ws.Range("K107").Value = 3814.2498
$ws.Range("L107").Value = 5085
$ws.Range("M107").Value = -1894.2498
$ws.Range("N107").Value = -8925
$ws.Range("H122").Value = 1597
$ws.Range("I122").Value = 1398.5
$ws.Range("K122").Value = 12586.5
$ws.Range("M122").Value = -10136.5
$ws.Range("H131").Value = 8078572
$ws.Range("I131").Value = 100000400
$ws.Range("K131").Value = 300001200
$ws.Range("M131").Value = -299996160
$ws.Range("H135").Value = 354.45
$ws.Range("I135").Value = 337.125
$ws.Range("J135").Value = 366
$ws.Range("K135").Value = 3034.125
$ws.Range("L135").Value = 3294
$ws.Range("M135").Value = -499.125
$ws.Range("N135").Value = -8364

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1626.1666
$ws.Range("I97").Value = 1160
$ws.Range("J97").Value = 2278.8
$ws.Range("K97").Value = 1160
$ws.Range("L97").Value = 2278.8
$ws.Range("M97").Value = -664
$ws.Range("N97").Value = -3270.8
$ws.Range("H102").Value = 4711.1665
$ws.Range("I102").Value = 4775.8184
$ws.Range("K102").Value = 4775.8184
$ws.Range("M102").Value = -3153.8184
$ws.Range("H126").Value = 1827332.5
$ws.Range("I126").Value = 2317889.5
$ws.Range("J126").Value = 145422.42
$ws.Range("K126").Value = 6953668.5
$ws.Range("L126").Value = 436267.26
$ws.Range("M126").Value = -6951198.5
$ws.Range("N126").Value = -441207.26
$ws.Range("H132").Value = 1925678
$ws.Range("I132").Value = 3498138.2
$ws.Range("J132").Value = 3782
$ws.Range("K132").Value = 10494414.6
$ws.Range("L132").Value = 11346
$ws.Range("M132").Value = -10491884.6
$ws.Range("N132").Value = -16406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3563.7222
$ws.Range("J7").Value = 4186.091
$ws.Range("L7").Value = 4186.091
$ws.Range("N7").Value = -4410.091
$ws.Range("H16").Value = 5305.2666
$ws.Range("I16").Value = 7267.9
$ws.Range("J16").Value = 1380
$ws.Range("K16").Value = 7267.9
$ws.Range("L16").Value = 1380
$ws.Range("M16").Value = -7097.9
$ws.Range("N16").Value = -1720
$ws.Range("H22").Value = 4285
$ws.Range("J22").Value = 5956
$ws.Range("L22").Value = 5956
$ws.Range("N22").Value = -6546
$ws.Range("H27").Value = 4285
$ws.Range("J27").Value = 5956
$ws.Range("L27").Value = 5956
$ws.Range("N27").Value = -6170
$ws.Range("H40").Value = 9129.299999999999
$ws.Range("I40").Value = 3659.2
$ws.Range("J40").Value = 14599.4
$ws.Range("K40").Value = 3659.2
$ws.Range("L40").Value = 14599.4
$ws.Range("M40").Value = -3523.2
$ws.Range("N40").Value = -14871.4
$ws.Range("H122").Value = 10988
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H126").Value = 3563.7222
$ws.Range("J126").Value = 4186.091
$ws.Range("L126").Value = 12558.273
$ws.Range("N126").Value = -17498.273
$ws.Range("H136").Value = 4267.5713
$ws.Range("I136").Value = 2609.923
$ws.Range("K136").Value = 7829.768999999999
$ws.Range("M136").Value = -5279.768999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4089.8333
$ws.Range("I126").Value = 3152.9473
$ws.Range("K126").Value = 9458.841899999999
$ws.Range("M126").Value = -6988.841899999999
$ws.Range("H136").Value = 13553681
$ws.Range("I136").Value = 23151500
$ws.Range("J136").Value = 3816.7646
$ws.Range("K136").Value = 69454500
$ws.Range("L136").Value = 11450.2938
$ws.Range("M136").Value = -69451950
$ws.Range("N136").Value = -16550.2938
